$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("output_and_resourceuse")

$ws.Range("F2").Value = 5.709401947836266
$ws.Range("G2").Value = 5.45
$ws.Range("H2").Value = 3.26
$ws.Range("I2").Value = 8.789999999999999

$ws.Range("I3").Value = 3.68

$ws.Range("F4").Value = 3.596301429414845
$ws.Range("G4").Value = 4.49
$ws.Range("H4").Value = 2.2
$ws.Range("I4").Value = 6.87

$ws.Range("I5").Value = 2.79

$ws.Range("I6").Value = 2.49

$ws.Range("H7").Value = 0.43
$ws.Range("I7").Value = 0.6

$ws.Range("D8").Value = 57214110.16
$ws.Range("G8").Value = 2.29
$ws.Range("H8").Value = 1.05
$ws.Range("I8").Value = 3.51
$ws.Range("J8").Value = 0.21
